$wb = $excel.ActiveWorkbook

# --- Sheet: Home win ---
$ws = $wb.Worksheets.Item('Home win')
$ws.Cells.Item(2,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'CHAMPIONSHIP'
$ws.Cells.Item(2,4).Value = 'Luton - Norwich'
$ws.Cells.Item(2,5).Value = 80
$ws.Cells.Item(2,6).Value = 2.25
$ws.Cells.Item(3,1).Value = '01-01-2025 13:00'
$ws.Cells.Item(3,2).Value = 'ENGLAND'
$ws.Cells.Item(3,3).Value = 'NATIONAL LEAGUE - NORTH'
$ws.Cells.Item(3,4).Value = 'Oxford City - Brackley Town'
$ws.Cells.Item(3,5).Value = 70
$ws.Cells.Item(3,6).Value = 2.45
$ws.Cells.Item(4,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(4,2).Value = 'ENGLAND'
$ws.Cells.Item(4,3).Value = 'NON LEAGUE PREMIER - ISTHMIAN'
$ws.Cells.Item(4,4).Value = 'Hashtag United - Billericay Town'
$ws.Cells.Item(4,5).Value = 73.3
$ws.Cells.Item(4,6).Value = 3.1
$ws.Cells.Item(5,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(5,2).Value = 'ENGLAND'
$ws.Cells.Item(5,3).Value = 'NON LEAGUE PREMIER - NORTHERN'
$ws.Cells.Item(5,4).Value = 'Lancaster City - Mickleover Sports'
$ws.Cells.Item(5,5).Value = 73.3
$ws.Cells.Item(5,6).Value = 1.75
$ws.Cells.Item(6,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(6,2).Value = 'ENGLAND'
$ws.Cells.Item(6,3).Value = 'NON LEAGUE PREMIER - NORTHERN'
$ws.Cells.Item(6,4).Value = 'Whitby Town - Stockton Town'
$ws.Cells.Item(6,5).Value = 73.3
$ws.Cells.Item(6,6).Value = 2.45
$ws.Cells.Item(7,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(7,2).Value = 'ENGLAND'
$ws.Cells.Item(7,3).Value = 'NON LEAGUE PREMIER - SOUTHERN CENTRAL'
$ws.Cells.Item(7,4).Value = 'Halesowen Town - Bromsgrove Sporting'
$ws.Cells.Item(7,5).Value = 80
$ws.Cells.Item(7,6).Value = 1.75

# --- Sheet: Draw ---
$ws = $wb.Worksheets.Item('Draw')
$ws.Cells.Item(2,1).Value = '31-12-2024 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'NATIONAL LEAGUE - NORTH'
$ws.Cells.Item(2,4).Value = 'Scunthorpe - King''s Lynn Town'
$ws.Cells.Item(2,5).Value = 60
$ws.Cells.Item(2,6).Value = 3.65
$ws.Cells.Item(3,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(3,2).Value = 'ENGLAND'
$ws.Cells.Item(3,3).Value = 'NATIONAL LEAGUE - NORTH'
$ws.Cells.Item(3,4).Value = 'Scarborough Athletic - Spennymoor Town'
$ws.Cells.Item(3,5).Value = 66.7
$ws.Cells.Item(3,6).Value = 3.5
$ws.Cells.Item(4,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(4,2).Value = 'ENGLAND'
$ws.Cells.Item(4,3).Value = 'NON LEAGUE PREMIER - ISTHMIAN'
$ws.Cells.Item(4,4).Value = 'Cray Wanderers - Dulwich Hamlet'
$ws.Cells.Item(4,5).Value = 66.7
$ws.Cells.Item(4,6).Value = 3.45
$ws.Cells.Item(5,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(5,2).Value = 'ENGLAND'
$ws.Cells.Item(5,3).Value = 'NON LEAGUE PREMIER - SOUTHERN CENTRAL'
$ws.Cells.Item(5,4).Value = 'Spalding United - St Ives Town'
$ws.Cells.Item(5,5).Value = 60
$ws.Cells.Item(5,6).Value = 3.45

# --- Sheet: Btts ---
$ws = $wb.Worksheets.Item('Btts')
$ws.Cells.Item(2,1).Value = '31-12-2024 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(2,4).Value = 'Solihull Moors - Tamworth'
$ws.Cells.Item(2,5).Value = 92
$ws.Cells.Item(2,6).Value = 1.75
$ws.Cells.Item(3,1).Value = '01-01-2025 17:30'
$ws.Cells.Item(3,2).Value = 'ENGLAND'
$ws.Cells.Item(3,3).Value = 'PREMIER LEAGUE'
$ws.Cells.Item(3,4).Value = 'Brentford - Arsenal'
$ws.Cells.Item(3,5).Value = 76.7
$ws.Cells.Item(3,6).Value = 1.75
$ws.Cells.Item(4,1).Value = '01-01-2025 17:30'
$ws.Cells.Item(4,2).Value = 'ENGLAND'
$ws.Cells.Item(4,3).Value = 'CHAMPIONSHIP'
$ws.Cells.Item(4,4).Value = 'Hull City - Middlesbrough'
$ws.Cells.Item(4,5).Value = 83.3
$ws.Cells.Item(4,6).Value = 1.73
$ws.Cells.Item(5,1).Value = '01-01-2025 20:00'
$ws.Cells.Item(5,2).Value = 'ENGLAND'
$ws.Cells.Item(5,3).Value = 'CHAMPIONSHIP'
$ws.Cells.Item(5,4).Value = 'Sunderland - Sheffield Utd'
$ws.Cells.Item(5,5).Value = 76.7
$ws.Cells.Item(5,6).Value = 2
$ws.Cells.Item(6,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(6,2).Value = 'ENGLAND'
$ws.Cells.Item(6,3).Value = 'LEAGUE TWO'
$ws.Cells.Item(6,4).Value = 'Port Vale - Cheltenham'
$ws.Cells.Item(6,5).Value = 80
$ws.Cells.Item(6,6).Value = 1.91
$ws.Cells.Item(7,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(7,2).Value = 'ENGLAND'
$ws.Cells.Item(7,3).Value = 'NATIONAL LEAGUE - SOUTH'
$ws.Cells.Item(7,4).Value = 'Weston-super-Mare - Chippenham Town'
$ws.Cells.Item(7,5).Value = 75
$ws.Cells.Item(7,6).Value = 1.75

# --- Sheet: Over_Under ---
$ws = $wb.Worksheets.Item('Over_Under')
$ws.Cells.Item(2,1).Value = '31-12-2024 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(2,4).Value = 'Solihull Moors - Tamworth'
$ws.Cells.Item(2,5).Value = 86.7
$ws.Cells.Item(2,6).Value = 1.7
$ws.Cells.Item(2,7).Value = 33.3
$ws.Cells.Item(2,8).Value = 2.75
$ws.Cells.Item(3,1).Value = '31-12-2024 08:00'
$ws.Cells.Item(3,2).Value = 'AUSTRALIA'
$ws.Cells.Item(3,3).Value = 'A-LEAGUE'
$ws.Cells.Item(3,4).Value = 'Central Coast Mariners - Melbourne City'
$ws.Cells.Item(3,5).Value = 80
$ws.Cells.Item(3,6).Value = 1.73
$ws.Cells.Item(3,7).Value = 50
$ws.Cells.Item(3,8).Value = 2.75
$ws.Cells.Item(4,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(4,2).Value = 'ENGLAND'
$ws.Cells.Item(4,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(4,4).Value = 'Rochdale - Altrincham'
$ws.Cells.Item(4,5).Value = 90
$ws.Cells.Item(4,6).Value = 1.85
$ws.Cells.Item(4,7).Value = 35
$ws.Cells.Item(4,8).Value = 3.25
$ws.Cells.Item(5,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(5,2).Value = 'ENGLAND'
$ws.Cells.Item(5,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(5,4).Value = 'York - Gateshead'
$ws.Cells.Item(5,5).Value = 80
$ws.Cells.Item(5,6).Value = 1.6
$ws.Cells.Item(5,7).Value = 60
$ws.Cells.Item(5,8).Value = 2.5
$ws.Cells.Item(6,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(6,2).Value = 'ENGLAND'
$ws.Cells.Item(6,3).Value = 'NATIONAL LEAGUE - NORTH'
$ws.Cells.Item(6,4).Value = 'Chorley - Buxton'
$ws.Cells.Item(6,5).Value = 85
$ws.Cells.Item(6,6).Value = 1.77
$ws.Cells.Item(6,7).Value = 50
$ws.Cells.Item(6,8).Value = 3
$ws.Cells.Item(7,1).Value = '01-01-2025 13:00'
$ws.Cells.Item(7,2).Value = 'ENGLAND'
$ws.Cells.Item(7,3).Value = 'NATIONAL LEAGUE - NORTH'
$ws.Cells.Item(7,4).Value = 'Oxford City - Brackley Town'
$ws.Cells.Item(7,5).Value = 86.7
$ws.Cells.Item(7,6).Value = 1.8
$ws.Cells.Item(7,7).Value = 40
$ws.Cells.Item(7,8).Value = 3.2
$ws.Cells.Item(8,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(8,2).Value = 'ENGLAND'
$ws.Cells.Item(8,3).Value = 'NATIONAL LEAGUE - NORTH'
$ws.Cells.Item(8,4).Value = 'South Shields - Darlington 1883'
$ws.Cells.Item(8,5).Value = 80
$ws.Cells.Item(8,6).Value = 1.7
$ws.Cells.Item(8,7).Value = 65
$ws.Cells.Item(8,8).Value = 2.88
$ws.Cells.Item(9,1).Value = '01-01-2025 15:00'
$ws.Cells.Item(9,2).Value = 'ENGLAND'
$ws.Cells.Item(9,3).Value = 'NON LEAGUE PREMIER - SOUTHERN CENTRAL'
$ws.Cells.Item(9,4).Value = 'Halesowen Town - Bromsgrove Sporting'
$ws.Cells.Item(9,5).Value = 80
$ws.Cells.Item(9,6).Value = 1.7
$ws.Cells.Item(9,7).Value = 35
$ws.Cells.Item(9,8).Value = ""

